$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.417.79'
$ws.Range("D3").Value = '1.944.11'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''242.67'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = '''0.612'
$ws.Range("E6").Value = '  -2.42%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").Value = '''57.52'
$ws.Range("E8").Value = '  -2.56%  '
$ws.Range("D9").Value = '''0.363'
$ws.Range("E9").Value = '  -3.43%  '
$ws.Range("E10").Value = '  +4.52%  '
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '2.229.14'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '''0.818'
$ws.Range("E13").Value = '  -5.44%  '
$ws.Range("D14").Value = '''21.29'
$ws.Range("E14").Value = '  -9.74%  '
$ws.Range("D15").Value = '''13.52'
$ws.Range("E15").Value = '  -3.50%  '
$ws.Range("D16").Value = '''5.19'
$ws.Range("E16").Value = '  -5.05%  '
$ws.Range("D17").Value = '1.942.23'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").Value = '36.356.06'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '0.0₃0874'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").Value = '''69.37'
$ws.Range("E20").Value = '  -1.83%  '
$ws.Range("D21").Value = '''228.76'
$ws.Range("E21").Value = '  -2.14%  '
$ws.Range("D22").Value = '''5.02'
$ws.Range("E22").Value = '  -6.05%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  -7.14%  '
$ws.Range("E25").Value = '  -1.02%  '
$ws.Range("E26").Value = '  -8.83%  '
$ws.Range("D27").Value = '''161.11'
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("E28").Value = '  +2.78%  '
$ws.Range("D29").Value = '''19.31'
$ws.Range("E29").Value = '  -2.77%  '
$ws.Range("E30").Value = '  -2.00%  '
$ws.Range("D31").Value = '''1.14'
$ws.Range("E31").Value = '  -6.44%  '
$ws.Range("D32").Value = '''4.63'
$ws.Range("E32").Value = '  -5.78%  '
$ws.Range("D33").Value = '''0.0635'
$ws.Range("E33").Value = '  +1.03%  '
$ws.Range("D34").Value = '''4.25'
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("E37").Value = '  +1.22%  '
$ws.Range("E38").Value = '  -4.92%  '
$ws.Range("D39").Value = '''3.07'
$ws.Range("E39").Value = '  +0.59%  '
$ws.Range("D40").Value = '''0.0976'
$ws.Range("E40").Value = '  +1.76%  '
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("E42").Value = '  -6.28%  '
$ws.Range("D43").Value = '''0.0211'
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("D44").Value = '''15.73'
$ws.Range("E44").Value = '  -2.99%  '
$ws.Range("D45").Value = '1.345.11'
$ws.Range("E45").Value = '  -2.92%  '
$ws.Range("E46").Value = '  -6.34%  '
$ws.Range("D47").Value = '''87.36'
$ws.Range("E47").Value = '  -5.79%  '
$ws.Range("D48").Value = '''7.15'
$ws.Range("E48").Value = '  -4.99%  '
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("D50").Value = '''44.43'
$ws.Range("E50").Value = '  -2.22%  '
$ws.Range("D51").Value = '2.119.95'
$ws.Range("E51").Value = '  -2.20%  '

# Reset style index for cells forced to text so no stray number-format/quote-prefix
# style is left applied (matches original un-styled inlineStr cells).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
